# Update product data rows for the "email a realizar pedido" (order request email) feature.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace product data ---
$ws.Range("A2").Value = "MB001002"
$ws.Range("B2").Value = "Peluche saco para dormir infan"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "50"
$ws.Range("F2").Value = "20"
$ws.Range("G2").Value = "19"
$ws.Range("H2").Value = "18"
$ws.Range("I2").Value = "17"
$ws.Range("J2").Value = "16"
$ws.Range("K2").ClearContents()
$ws.Range("M2:T2").ClearContents()
$ws.Range("AB2:AE2").ClearContents()

# --- Row 3: replace product data ---
$ws.Range("A3").Value = "LapHP"
$ws.Range("B3").Value = "Laptop HP"
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = "75"
$ws.Range("F3").Value = "1000"
$ws.Range("G3").Value = "1250.3"
$ws.Range("H3").Value = "1180.99"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").ClearContents()
$ws.Range("M3:T3").ClearContents()
$ws.Range("AB3:AE3").ClearContents()

# --- Row 4: removed entirely (was DISPENSADOR DE JUGO X3) ---
$ws.Rows.Item(4).Delete()
